$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2 held the latitude as a text string ("39.296717743018434, "); replace it
# with the real numeric latitude value.
$ws.Range("C2").Value = 39.296717743018398

# Column C (Latitude) is now numeric-sized; apply the best-fit width Excel
# computed for that content (stored as ~11 characters wide).
$ws.Columns.Item(3).ColumnWidth = 10.17

# The sheet view's last-known selection moved to H5.
$ws.Range("H5").Select()
